# Release Log Form - F6.xlsx : "feat: sops Update 4"
#
# 1. Rename sheet "F-SW-CR-06" -> "F-SW-SD-06"
# 2. Un-hide "Sheet2"
# 3. Repoint the Print_Area defined name at the renamed sheet
# 4. Move the active-cell selection on sheet 1 from J2 -> E24
# 5. Update the footer date stamp "0/0/2025" -> "01/10/2025"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename the first worksheet.
$ws1.Name = "F-SW-SD-06"

# 2. Make Sheet2 visible again.
$ws2.Visible = -1

# 3. Update the workbook-level Print_Area defined name so it still points at
#    the (renamed) first sheet.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "='F-SW-SD-06'!`$A`$1:`$G`$31"
    }
}

# 4. Update the selection / active cell on sheet 1.
$ws1.Activate()
$ws1.Range("E24").Select()

# 5. Update the right-hand footer text with the new revision date.
$ws1.PageSetup.RightFooter = "&14Rev:0(01/10/2025)"
